$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34
$ws.Range("D34").Value = 200
$ws.Range("E34").Value = "North"
$ws.Range("F34").Value = 0

# Row 35
$ws.Range("D35").Value = 800
$ws.Range("E35").Value = "East"
$ws.Range("F35").Value = 90

# Row 36
$ws.Range("D36").Value = 1400
$ws.Range("E36").Value = "South"
$ws.Range("F36").Value = 180

# Row 37
$ws.Range("D37").Value = 2000
$ws.Range("E37").Value = "west"
$ws.Range("F37").Value = 270

# Row 38
$ws.Range("D38").Value = 2600
$ws.Range("E38").Value = "North"
$ws.Range("F38").Value = 360

# Row 33: header "Millivolts" in D33 (added last so shared-string index matches)
$ws.Range("D33").Value = "Millivolts"

# Update view: selection matches the new data block
$ws.Range("D33:E38").Select() | Out-Null
